# Add support for a second fuel / second efficiency on PowerPlants.
#
# 1) Sheet "PowerPlants" (sheet3): insert a new "fuel2" [text] column
#    right after the existing "fuel" column.
# 2) Sheet "PowerPlantsPerformance" (sheet4): insert a new "Efficiency2"
#    [%] column right after the existing "Efficiency" column.
# 3) Move the active worksheet / selection from PowerPlants to
#    PowerPlantsPerformance, as left by the author after editing.

$wb = $excel.ActiveWorkbook

# --- PowerPlants sheet: insert "fuel2" column after "fuel" (column C) ---
$wsPowerPlants = $wb.Worksheets.Item("PowerPlants")
$wsPowerPlants.Columns.Item(3).Insert()
$wsPowerPlants.Cells.Item(1, 3).Value = "fuel2"
$wsPowerPlants.Cells.Item(2, 3).Value = "[text]"
$wsPowerPlants.Columns.Item(3).ColumnWidth = 14.42

# --- PowerPlantsPerformance sheet: insert "Efficiency2" column after "Efficiency" (column D -> new column E) ---
$wsPerformance = $wb.Worksheets.Item("PowerPlantsPerformance")
$wsPerformance.Columns.Item(5).Insert()
$wsPerformance.Cells.Item(1, 5).Value = "Efficiency2"
$wsPerformance.Cells.Item(2, 5).Value = "[%]"
$wsPerformance.Columns.Item(5).ColumnWidth = 11.9

# --- Update selections left on PowerPlants and move focus to PowerPlantsPerformance ---
$wsPowerPlants.Activate()
$wsPowerPlants.Range("C3").Select()

$wsPerformance.Activate()
$wsPerformance.Range("E2").Select()
